$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 21:52"

# Update Canary Islands provinces rows with new "Casos activos" (C) and "Recuperados" (D) values
$rows = @(32, 47, 56, 57, 59, 62, 64)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 249
    $ws.Cells.Item($r, 4).Value = 1422
}
